# Recalculated profit figures (columns H-N) across the Leve tables on every
# sheet, reflecting refreshed market-board prices from the scheduled data pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 840.4545000000001
$ws.Range("I11").Value = 840.4545000000001
$ws.Range("K11").Value = 840.4545000000001
$ws.Range("M11").Value = -700.4545000000001
$ws.Range("H17").Value = 3036.923
$ws.Range("J17").Value = 3036.923
$ws.Range("L17").Value = 9110.769
$ws.Range("N17").Value = -9446.769
$ws.Range("H28").Value = 176.15384
$ws.Range("I28").Value = 157.5
$ws.Range("K28").Value = 157.5
$ws.Range("M28").Value = 327.5
$ws.Range("H48").Value = 2760.2856
$ws.Range("I48").Value = 487
$ws.Range("J48").Value = 3669.6
$ws.Range("K48").Value = 1461
$ws.Range("L48").Value = 11008.8
$ws.Range("M48").Value = -1169
$ws.Range("N48").Value = -11592.8
$ws.Range("H56").Value = 2760.2856
$ws.Range("I56").Value = 487
$ws.Range("J56").Value = 3669.6
$ws.Range("K56").Value = 1461
$ws.Range("L56").Value = 11008.8
$ws.Range("M56").Value = -927
$ws.Range("N56").Value = -12076.8
$ws.Range("H58").Value = 819.4286
$ws.Range("I58").Value = 182.14285
$ws.Range("J58").Value = 1456.7142
$ws.Range("K58").Value = 546.4285500000001
$ws.Range("L58").Value = 4370.142599999999
$ws.Range("M58").Value = -396.4285500000001
$ws.Range("N58").Value = -4670.142599999999
$ws.Range("H59").Value = 444
$ws.Range("I59").Value = 444
$ws.Range("K59").Value = 1332
$ws.Range("M59").Value = -775
$ws.Range("H74").Value = 10697.167
$ws.Range("I74").Value = 8685.571
$ws.Range("J74").Value = 11977.272
$ws.Range("K74").Value = 8685.571
$ws.Range("L74").Value = 11977.272
$ws.Range("M74").Value = -7749.571
$ws.Range("N74").Value = -13849.272
$ws.Range("H77").Value = 10697.167
$ws.Range("I77").Value = 8685.571
$ws.Range("J77").Value = 11977.272
$ws.Range("K77").Value = 43427.855
$ws.Range("L77").Value = 59886.36
$ws.Range("M77").Value = -38747.855
$ws.Range("N77").Value = -69246.36
$ws.Range("H92").Value = 83.90000000000001
$ws.Range("I92").Value = 83.90000000000001
$ws.Range("K92").Value = 83.90000000000001
$ws.Range("M92").Value = 1164.1
$ws.Range("H98").Value = 1388.5385
$ws.Range("I98").Value = 1388.5385
$ws.Range("K98").Value = 1388.5385
$ws.Range("M98").Value = 109.4614999999999
$ws.Range("H107").Value = 515.75
$ws.Range("I107").Value = 419.1
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 419.1
$ws.Range("L107").Value = 999
$ws.Range("M107").Value = 1500.9
$ws.Range("N107").Value = -4839
$ws.Range("H122").Value = 1388.5385
$ws.Range("I122").Value = 1388.5385
$ws.Range("K122").Value = 4165.6155
$ws.Range("M122").Value = -1715.6155
$ws.Range("H125").Value = 1079415.4
$ws.Range("I125").Value = 1847287.9
$ws.Range("K125").Value = 16625591.1
$ws.Range("M125").Value = -16623131.1
$ws.Range("H132").Value = 2310.4736
$ws.Range("I132").Value = 2185.6428
$ws.Range("K132").Value = 6556.928400000001
$ws.Range("M132").Value = -4026.928400000001
$ws.Range("H138").Value = 3350.4194
$ws.Range("I138").Value = 2723.8333
$ws.Range("J138").Value = 3500.8
$ws.Range("K138").Value = 8171.499899999999
$ws.Range("L138").Value = 10502.4
$ws.Range("M138").Value = -3031.499899999999
$ws.Range("N138").Value = -20782.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H32").Value = 10348.417
$ws.Range("I32").Value = 8445.143
$ws.Range("J32").Value = 23671.334
$ws.Range("K32").Value = 8445.143
$ws.Range("L32").Value = 23671.334
$ws.Range("M32").Value = -8158.143
$ws.Range("N32").Value = -24245.334
$ws.Range("H74").Value = 3186.3125
$ws.Range("I74").Value = 2541
$ws.Range("J74").Value = 6671
$ws.Range("K74").Value = 2541
$ws.Range("L74").Value = 6671
$ws.Range("M74").Value = -1667
$ws.Range("N74").Value = -8419
$ws.Range("H77").Value = 3186.3125
$ws.Range("I77").Value = 2541
$ws.Range("J77").Value = 6671
$ws.Range("K77").Value = 12705
$ws.Range("L77").Value = 33355
$ws.Range("M77").Value = -8337
$ws.Range("N77").Value = -42091
$ws.Range("H132").Value = 1895.814
$ws.Range("I132").Value = 1718
$ws.Range("J132").Value = 4266.6665
$ws.Range("K132").Value = 5154
$ws.Range("L132").Value = 12799.9995
$ws.Range("M132").Value = -2624
$ws.Range("N132").Value = -17859.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 199950
$ws.Range("J43").Value = 199950
$ws.Range("L43").Value = 199950
$ws.Range("N43").Value = -200312
$ws.Range("H86").Value = 1453.1034
$ws.Range("I86").Value = 1289.25
$ws.Range("J86").Value = 2239.6
$ws.Range("K86").Value = 1289.25
$ws.Range("L86").Value = 2239.6
$ws.Range("M86").Value = -166.25
$ws.Range("N86").Value = -4485.6
$ws.Range("H89").Value = 1453.1034
$ws.Range("I89").Value = 1289.25
$ws.Range("J89").Value = 2239.6
$ws.Range("K89").Value = 6446.25
$ws.Range("L89").Value = 11198
$ws.Range("M89").Value = -830.25
$ws.Range("N89").Value = -22430
$ws.Range("H134").Value = 3927.3333
$ws.Range("I134").Value = 3908.1929
$ws.Range("K134").Value = 11724.5787
$ws.Range("M134").Value = -9189.5787

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5735.0835
$ws.Range("I31").Value = 5468.5
$ws.Range("J31").Value = 5823.9443
$ws.Range("K31").Value = 5468.5
$ws.Range("L31").Value = 5823.9443
$ws.Range("M31").Value = -5173.5
$ws.Range("N31").Value = -6413.9443
$ws.Range("H34").Value = 5735.0835
$ws.Range("I34").Value = 5468.5
$ws.Range("J34").Value = 5823.9443
$ws.Range("K34").Value = 5468.5
$ws.Range("L34").Value = 5823.9443
$ws.Range("M34").Value = -5266.5
$ws.Range("N34").Value = -6227.9443
$ws.Range("H68").Value = 30637.5
$ws.Range("J68").Value = 27375
$ws.Range("L68").Value = 27375
$ws.Range("N68").Value = -28873
$ws.Range("H71").Value = 30637.5
$ws.Range("J71").Value = 27375
$ws.Range("L71").Value = 82125
$ws.Range("N71").Value = -89613
$ws.Range("H74").Value = 37493.332
$ws.Range("J74").Value = 37493.332
$ws.Range("L74").Value = 37493.332
$ws.Range("N74").Value = -39241.332
$ws.Range("H77").Value = 37493.332
$ws.Range("J77").Value = 37493.332
$ws.Range("L77").Value = 112479.996
$ws.Range("N77").Value = -121215.996
$ws.Range("H134").Value = 11250
$ws.Range("I134").Value = 7500
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 22500
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -19965
$ws.Range("N134").Value = -50070
$ws.Range("H141").Value = 249011
$ws.Range("J141").Value = 249011
$ws.Range("L141").Value = 249011
$ws.Range("N141").Value = -259371

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 14823893
$ws.Range("I4").Value = 10954658
$ws.Range("J4").Value = 34170064
$ws.Range("K4").Value = 32863974
$ws.Range("L4").Value = 102510192
$ws.Range("M4").Value = -32863862
$ws.Range("N4").Value = -102510416
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H54").Value = 499.33334
$ws.Range("I54").Value = 499.5
$ws.Range("J54").Value = 499
$ws.Range("K54").Value = 1498.5
$ws.Range("L54").Value = 1497
$ws.Range("M54").Value = -939.5
$ws.Range("N54").Value = -2615
$ws.Range("H55").Value = 988.1667
$ws.Range("I55").Value = 980
$ws.Range("J55").Value = 996.3333
$ws.Range("K55").Value = 2940
$ws.Range("L55").Value = 2988.9999
$ws.Range("M55").Value = -2763
$ws.Range("N55").Value = -3342.9999
$ws.Range("H114").Value = 5000
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 5000
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 15000
$ws.Range("M114").ClearContents()
$ws.Range("N114").Value = -21508

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1917.3043
$ws.Range("I122").Value = 1829.9
$ws.Range("K122").Value = 5489.700000000001
$ws.Range("M122").Value = -3039.700000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1450.6666
$ws.Range("J16").Value = 894.75
$ws.Range("L16").Value = 894.75
$ws.Range("N16").Value = -1234.75
$ws.Range("H22").Value = 3342.7144
$ws.Range("I22").Value = 3549.5
$ws.Range("K22").Value = 3549.5
$ws.Range("M22").Value = -3254.5
$ws.Range("H27").Value = 3342.7144
$ws.Range("I27").Value = 3549.5
$ws.Range("K27").Value = 3549.5
$ws.Range("M27").Value = -3442.5
$ws.Range("H43").Value = 44979
$ws.Range("J43").Value = 44979
$ws.Range("L43").Value = 44979
$ws.Range("N43").Value = -45365
$ws.Range("H74").Value = 45998
$ws.Range("J74").Value = 75000
$ws.Range("L74").Value = 75000
$ws.Range("N74").Value = -76996
$ws.Range("H77").Value = 45998
$ws.Range("J77").Value = 75000
$ws.Range("L77").Value = 225000
$ws.Range("N77").Value = -234984

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1076.8889
$ws.Range("I100").Value = 977.4286
$ws.Range("K100").Value = 1954.8572
$ws.Range("M100").Value = -1413.8572
$ws.Range("H132").Value = 7245.3076
$ws.Range("I132").Value = 6364.2607
$ws.Range("K132").Value = 19092.7821
$ws.Range("M132").Value = -16562.7821
$ws.Range("H136").Value = 5343.64
$ws.Range("I136").Value = 5088.6113
$ws.Range("J136").Value = 5999.4287
$ws.Range("K136").Value = 15265.8339
$ws.Range("L136").Value = 17998.2861
$ws.Range("M136").Value = -12715.8339
$ws.Range("N136").Value = -23098.2861
